$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45189 = 2023-09-20).
# Update every data row (2 through 358) to the new date serial 45190 (2023-09-21).
$ws.Range("C2:C358").Value = 45190
